# iCub3 skin GUI cleanup files for upperarms
# Update the reference-geometry parameters (O5/O7/O9) on the
# "left upperarm front" sheet, which ripple through the dependent
# formulas on that sheet and on the "final ini" sheet, and restore
# the two sheets' last-used selections.

$wb = $excel.ActiveWorkbook

# --- "left upperarm front" -------------------------------------------------
$front = $wb.Worksheets.Item("left upperarm front")

$front.Range("O5").Value = 180
$front.Range("O7").Value = 130
$front.Range("O9").Value = 190

[void]$front.Range("O8").Select()

# --- "final ini" -------------------------------------------------------------
$finalIni = $wb.Worksheets.Item("final ini")

[void]$finalIni.Range("A13:G20").Select()
